# Update gh-pages to output generated at 456a3b4
#
# This script re-applies the numeric "want-to-go" count refresh (column F)
# across the four sheets, plus the brand-new 演出 (Performance) row that was
# scraped in after row 21.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: set a cell's value while suppressing Excel's "looks like a date"
# auto-conversion for plain-text columns (B date-label text, E time-range
# text, H links, I image paths). We briefly force Text format, write the
# literal string, then clear formatting back to the sheet default so the
# exported cell keeps the original (unstyled) look.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions) — column F ("想去人数") refresh
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 185
    6  = 159
    8  = 758
    9  = 4135
    11 = 53
    14 = 5908
    15 = 460
    16 = 2290
    18 = 159
    19 = 447
    20 = 8892
    22 = 1911
    25 = 2377
    26 = 1374
    27 = 221
    28 = 1924
    31 = 320
    41 = 88
    43 = 1496
    44 = 2405
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances) — column F refresh on row 20, plus a brand
# new row 21 entry (displayed index 21) appended after the existing data.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F20").Value = 4

$newRow = 22
$ws2.Range("A21").Copy()
$ws2.Range("A$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A$newRow").Value = 21
Set-TextValue $ws2.Range("B$newRow") "2024-11-30"
Set-TextValue $ws2.Range("C$newRow") "北京·花たん 2024 LIVE in Beijing"
Set-TextValue $ws2.Range("D$newRow") "复兴路69号院2号136、G23室 Mao Livehouse北京五棵松店"
Set-TextValue $ws2.Range("E$newRow") "2024.11.30 14:00-11.30 15:30"
$ws2.Range("F$newRow").Value = 0
$ws2.Range("G$newRow").Value = 380
Set-TextValue $ws2.Range("H$newRow") "https://show.bilibili.com/platform/detail.html?id=90341"
Set-TextValue $ws2.Range("I$newRow") "//i2.hdslb.com/bfs/openplatform/202408/wfGEn3sY1722910561352.jpeg"

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) — column F refresh
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 679
$ws3.Range("F3").Value = 876

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) — column F refresh (same source counts as
# above, merged across every category, plus the new row's 0 is appended
# within the 演出 section but isn't part of this aggregate sheet's diff).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 679
    4  = 876
    8  = 159
    12 = 758
    13 = 4135
    14 = 4135
    15 = 53
    18 = 5908
    19 = 460
    20 = 2290
    22 = 159
    23 = 447
    24 = 8892
    27 = 1911
    29 = 2377
    30 = 1374
    31 = 221
    32 = 1924
    35 = 320
    41 = 88
    43 = 1497
    44 = 2405
    48 = 4
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
